$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 383; this shifts the previous rows
# 383..465 down to 385..467 (matching the dimension growing from
# A1:R465 to A1:R467).
$ws.Rows("383:384").Insert()

# New row 383: a fresh Cebollín record for Provincia de Cautín.
$ws.Cells.Item(383, 1).Value = 10
$ws.Cells.Item(383, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(383, 3).Value = "La Araucanía"
$ws.Cells.Item(383, 4).Value = 44855
$ws.Cells.Item(383, 5).Value = 9
$ws.Cells.Item(383, 6).Value = 100112037
$ws.Cells.Item(383, 7).Value = "Cebollín"
$ws.Cells.Item(383, 8).Value = "Sin especificar"
$ws.Cells.Item(383, 9).Value = "Primera"
$ws.Cells.Item(383, 10).Value = 50
$ws.Cells.Item(383, 11).Value = 8000
$ws.Cells.Item(383, 12).Value = 9000
$ws.Cells.Item(383, 13).Value = 8600
$ws.Cells.Item(383, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(383, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(383, 16).Value = 717
$ws.Cells.Item(383, 17).Value = 12
$ws.Cells.Item(383, 18).Value = "Hortaliza"

# New row 384: a fresh Cebollín record for Región Metropolitana.
$ws.Cells.Item(384, 1).Value = 10
$ws.Cells.Item(384, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(384, 3).Value = "La Araucanía"
$ws.Cells.Item(384, 4).Value = 44855
$ws.Cells.Item(384, 5).Value = 9
$ws.Cells.Item(384, 6).Value = 100112037
$ws.Cells.Item(384, 7).Value = "Cebollín"
$ws.Cells.Item(384, 8).Value = "Sin especificar"
$ws.Cells.Item(384, 9).Value = "Primera"
$ws.Cells.Item(384, 10).Value = 70
$ws.Cells.Item(384, 11).Value = 7000
$ws.Cells.Item(384, 12).Value = 8000
$ws.Cells.Item(384, 13).Value = 7429
$ws.Cells.Item(384, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(384, 15).Value = "Región Metropolitana"
$ws.Cells.Item(384, 16).Value = 619
$ws.Cells.Item(384, 17).Value = 12
$ws.Cells.Item(384, 18).Value = "Hortaliza"
